$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 39: titration/CRM accuracy data for the 2021-06-01 run (TP2 run 6 / run 4)
$ws.Range("A39").Value = 20210601
$ws.Range("B39").Value = 2222.2813428003801
$ws.Range("C39").Value = 2224.4699999999998
$ws.Range("D39").Formula = "=100*(B39-C39)/C39"
$ws.Range("E39").Value = 180
$ws.Range("F39").Value = "CRM opened 20210526"

# Leave the selection where the author left it after entering the new row
[void]$ws.Range("F40").Select()
